{"js": "// Load the body paragraphs so we can address them by index.\nconst body = context.document.body;\nbody.paragraphs.load(\"items,text\");\nawait context.sync();\n\nconst paras = body.paragraphs.items;\n\n// 1) \"get it digesting changes to the documents\" -> \"include langchain functionality\"\nparas[0].insertText(\"include langchain functionality\", \"Replace\");\n\n// 2) \"include references\" -> \"process sharepoint documents\"\nparas[1].insertText(\"process sharepoint documents\", \"Replace\");\n\n// 3) \"stateful persistence\" -> \"deploy using docker\"\nparas[2].insertText(\"deploy using docker\", \"Replace\");\n\n// 4) \"include langchain functionality\" (3 runs incl. proofErr) -> \"document the code\"\nparas[3].insertText(\"document the code\", \"Replace\");\n\n// 5) \"process sharepoint documents\" (3 runs incl. proofErr) ->\n//    \"Monitor performance, cost, and user satisfaction during testing.\"\nparas[4].insertText(\"Monitor performance, cost, and user satisfaction during testing.\", \"Replace\");\n\n// 6) \"understand how our documents are not being used to train the LLMs\" ->\n//    \"Do the embedding process separately, not every time the web page loads\"\nparas[5].insertText(\"Do the embedding process separately, not every time the web page loads\", \"Replace\");\n\n// 7-10) Remove the trailing paragraphs that are no longer part of the list:\n//    \"deploy using docker\", \"document the code\", \"show Chris, Geoff, Tim, Siobhan\",\n//    and the trailing empty list-styled paragraph.\nparas[6].delete();\nparas[7].delete();\nparas[8].delete();\nparas[9].delete();\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\nfunction Set-ParaText($para, $newText) {\n    # Scope Find/Replace to this paragraph's own Range so the whole-document\n    # text doesn't accidentally get matched (some of the new strings duplicate\n    # text that other paragraphs are changed to/from). This also correctly\n    # collapses a multi-run paragraph (e.g. one split by a spell-check\n    # <w:proofErr/> run) down to the replacement text.\n    $rng = $para.Range\n    $f = $rng.Find\n    $f.ClearFormatting()\n    $f.Replacement.ClearFormatting()\n    # wdFindContinue = 1, wdReplaceAll = 2\n    $f.Execute($para.Range.Text.TrimEnd([char]13, [char]7), $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n\n# 1) \"get it digesting changes to the documents\" -> \"include langchain functionality\"\nSet-ParaText $d.Paragraphs.Item(1) \"include langchain functionality\"\n\n# 2) \"include references\" -> \"process sharepoint documents\"\nSet-ParaText $d.Paragraphs.Item(2) \"process sharepoint documents\"\n\n# 3) \"stateful persistence\" -> \"deploy using docker\"\nSet-ParaText $d.Paragraphs.Item(3) \"deploy using docker\"\n\n# 4) \"include langchain functionality\" (3 runs incl. proofErr) -> \"document the code\"\nSet-ParaText $d.Paragraphs.Item(4) \"document the code\"\n\n# 5) \"process sharepoint documents\" (3 runs incl. proofErr) ->\n#    \"Monitor performance, cost, and user satisfaction during testing.\"\nSet-ParaText $d.Paragraphs.Item(5) \"Monitor performance, cost, and user satisfaction during testing.\"\n\n# 6) \"understand how our documents are not being used to train the LLMs\" ->\n#    \"Do the embedding process separately, not every time the web page loads\"\nSet-ParaText $d.Paragraphs.Item(6) \"Do the embedding process separately, not every time the web page loads\"\n\n# 7-10) Remove the trailing paragraphs that are no longer part of the list:\n#    \"deploy using docker\", \"document the code\", \"show Chris, Geoff, Tim, Siobhan\",\n#    and the trailing empty list-styled paragraph. Delete from the end backward\n#    so earlier indices stay valid.\n$d.Paragraphs.Item(10).Range.Delete()\n$d.Paragraphs.Item(9).Range.Delete()\n$d.Paragraphs.Item(8).Range.Delete()\n$d.Paragraphs.Item(7).Range.Delete()\n"}
